$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'229.73"
$ws.Range("D3").Value = "'22.30"
$ws.Range("D4").Value = "'5.241"
$ws.Range("D5").Value = "'0.05555"
$ws.Range("D6").Value = "'3.379"
$ws.Range("D7").Value = "'6.468"
$ws.Range("D8").Value = "'1.059"
$ws.Range("D9").Value = "'0.7801"
$ws.Range("D10").Value = "'0.1372"
$ws.Range("D11").Value = "'0.07325"
$ws.Range("D12").Value = "'0.03158"
$ws.Range("D13").Value = "'0.02944"
$ws.Range("D14").Value = "'0.09261"
$ws.Range("D15").Value = "'0.001659"
$ws.Range("D16").Value = "'3.274"
$ws.Range("D17").Value = "'0.04785"
$ws.Range("D18").Value = "'0.0005898"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006220"
$ws.Range("D20").Value = "'0.005236"
$ws.Range("D23").Value = "'3.910"
$ws.Range("D26").Value = "'0.1243"
$ws.Range("D27").Value = "'0.0004999"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Value = "'0.03991"
$ws.Range("D41").Value = "'0.007152"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.003499"
$ws.Range("D44").Value = "'0.01005"
$ws.Range("D45").Value = "'0.00005437"
$ws.Range("D47").Value = "'0.7850"
$ws.Range("D48").Value = "'0.05172"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
